# Updated cryptos list on Fri Jun  2 03:53:48 UTC 2023 with GitHub Actions
#
# Column D ("Price") cells hold plain text (not numbers) in this workbook,
# so many values ("26.932.32", "1.000", "0.000008486", trailing-zero
# figures like "0.3720", ...) would be silently re-interpreted/garbled if
# assigned while the cell is in its default "General" number format.
# Temporarily switching the cell to Text ("@") before the assignment keeps
# the literal string, then resetting the style back to "Normal" removes
# the now-unneeded text-format style so the cell's style index is left
# exactly as it was (unstyled, same as every other data cell in the row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $value) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

# Rows 13 and 14 swap coin identity (Litecoin <-> WrappedEther) and get
# refreshed price/volume figures.
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextValue "D13" "1.885.30"
$ws.Range("E13").Value = "  +1.54%  "

$ws.Range("B14").Value = "Litecoin"
$ws.Range("C14").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
Set-TextValue "D14" "94.97"
$ws.Range("E14").Value = "  +4.88%  "

# Refreshed price (column D) and volume (1h) percentage (column E) figures
# for the remaining rows.
Set-TextValue "D2" "26.932.32"
$ws.Range("E2").Value = "  +0.49%  "

Set-TextValue "D3" "1.875.66"
$ws.Range("E3").Value = "  +1.03%  "

Set-TextValue "D4" "0.9998"

Set-TextValue "D5" "306.19"
$ws.Range("E5").Value = "  +0.36%  "

Set-TextValue "D6" "0.9994"
$ws.Range("E6").Value = "  -0.11%  "

Set-TextValue "D7" "0.5159"
$ws.Range("E7").Value = "  +1.49%  "

Set-TextValue "D8" "0.3720"
$ws.Range("E8").Value = "  +1.90%  "

$ws.Range("E9").Value = "  +0.91%  "

Set-TextValue "D10" "0.8982"
$ws.Range("E10").Value = "  +1.25%  "

$ws.Range("E11").Value = "  -0.15%  "

Set-TextValue "D12" "0.07567"
$ws.Range("E12").Value = "  +0.91%  "

Set-TextValue "D15" "5.248"
$ws.Range("E15").Value = "  +0.21%  "

Set-TextValue "D16" "1.000"
$ws.Range("E16").Value = "  -0.11%  "

Set-TextValue "D17" "0.000008486"
$ws.Range("E17").Value = "  -0.29%  "

$ws.Range("E18").Value = "  +1.50%  "

$ws.Range("E19").Value = "  -0.09%  "

Set-TextValue "D20" "26.964.57"
$ws.Range("E20").Value = "  +0.45%  "

Set-TextValue "D21" "5.029"
$ws.Range("E21").Value = "  +0.56%  "

Set-TextValue "D22" "2.085.00"
$ws.Range("E22").Value = "  -0.64%  "

$ws.Range("E23").Value = "  +1.25%  "

Set-TextValue "D24" "6.434"
$ws.Range("E24").Value = "  -0.02%  "

Set-TextValue "D25" "145.96"

Set-TextValue "D26" "1.781"
$ws.Range("E26").Value = "  -2.04%  "

Set-TextValue "D27" "18.04"
$ws.Range("E27").Value = "  +1.19%  "

Set-TextValue "D28" "2.109"
$ws.Range("E28").Value = "  +3.28%  "

Set-TextValue "D29" "114.54"
$ws.Range("E29").Value = "  +1.57%  "

Set-TextValue "D30" "4.895"
$ws.Range("E30").Value = "  +5.01%  "

Set-TextValue "D31" "4.745"
$ws.Range("E31").Value = "  +2.78%  "

Set-TextValue "D32" "0.09180"
$ws.Range("E32").Value = "  -0.56%  "

Set-TextValue "D33" "0.05030"
$ws.Range("E33").Value = "  -1.40%  "

Set-TextValue "D34" "0.7534"
$ws.Range("E34").Value = "  +3.04%  "

Set-TextValue "D35" "2.994"
$ws.Range("E35").Value = "  -2.35%  "

$ws.Range("E36").Value = "  +1.94%  "

Set-TextValue "D37" "3.287"
$ws.Range("E37").Value = "  +3.11%  "

Set-TextValue "D38" "0.01991"
$ws.Range("E38").Value = "  -1.06%  "

Set-TextValue "D39" "0.5575"
$ws.Range("E39").Value = "  +5.51%  "

Set-TextValue "D40" "2.478"
$ws.Range("E40").Value = "  +0.75%  "

Set-TextValue "D41" "1.073"
$ws.Range("E41").Value = "  +0.08%  "

Set-TextValue "D42" "6.568"
$ws.Range("E42").Value = "  +1.82%  "

Set-TextValue "D43" "8.759"
$ws.Range("E43").Value = "  +5.02%  "

Set-TextValue "D44" "116.12"
$ws.Range("E44").Value = "  -1.43%  "

Set-TextValue "D45" "0.1501"
$ws.Range("E45").Value = "  +2.08%  "

Set-TextValue "D46" "0.4765"
$ws.Range("E46").Value = "  +2.74%  "

Set-TextValue "D47" "0.9989"
$ws.Range("E47").Value = "  -0.13%  "

Set-TextValue "D48" "10.08"
$ws.Range("E48").Value = "  +1.69%  "

Set-TextValue "D49" "1.562"
$ws.Range("E49").Value = "  +0.64%  "

Set-TextValue "D50" "37.15"
$ws.Range("E50").Value = "  +0.38%  "

Set-TextValue "D51" "63.33"
$ws.Range("E51").Value = "  +0.47%  "
